# Understanding Git.pptx - "Reword explanation of commit parents."
#
# 1) Handout master & notes master: the cached "datetimeFigureOut" date
#    field text moves from 4/13/2017 to 4/18/2017.
# 2) Slide 11 ("Commits in git"): reword the first three bullets that
#    explain commit parents.

$p = $ppt.ActivePresentation

# --- 1a) Handout Master date placeholder -----------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "4/18/2017"

# --- 1b) Notes Master date placeholder --------------------------------
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "4/18/2017"

# --- 2) Slide 11 content placeholder bullets --------------------------
$s = $p.Slides.Item(11)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$tr.Paragraphs(1, 1).Text = "Commits have parents."
$tr.Paragraphs(2, 1).Text = "A parent is just the snapshot(s) prior to current one."
$tr.Paragraphs(3, 1).Text = "A commit with multiple parents is a merge commit. It combines the changes in both of its parents."
